# Fruta / hortaliza, semanal
# A new weekly record is inserted as row 74 ("Terminal Hortofrutícola Agro
# Chillán" - Poroto granado), pushing the existing rows 74-87 down to 75-88.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, shifting rows 74-87 down to 75-88.
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new weekly entry.
$ws.Range("A74").Value2 = 7
$ws.Range("B74").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C74").Value2 = "Ñuble"
$ws.Range("D74").Value2 = 44617
$ws.Range("E74").Value2 = 16
$ws.Range("F74").Value2 = 100112030
$ws.Range("G74").Value2 = "Poroto granado"
$ws.Range("H74").Value2 = "Sin especificar"
$ws.Range("I74").Value2 = "Primera"
$ws.Range("J74").Value2 = 120
$ws.Range("K74").Value2 = 20000
$ws.Range("L74").Value2 = 20000
$ws.Range("M74").Value2 = 20000
$ws.Range("N74").Value2 = "$/saco 25 kilos"
$ws.Range("O74").Value2 = "Provincia de Diguillín"
$ws.Range("P74").Value2 = 800
$ws.Range("Q74").Value2 = 25
$ws.Range("R74").Value2 = "Hortaliza"
